# "made progress in creating temp data table"
#
# 1. Re-wrap the instructions text in cell A26 of the "REPORT" sheet in
#    literal double quotes.
# 2. Move the view/selection: scroll so row 13 is the top-left visible
#    row, and select A28:H31 (active cell A28) instead of A26:I26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPORT")
$ws.Activate()

$oldText = $ws.Range("A26").Value2
$ws.Range("A26").Value = '"' + $oldText + '"'

$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("A28:H31").Select()
